$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (3-45) form a cyclic rotation: each row's content moves down
# by one row, and the last row (45) wraps around to become the new row 3.
# i.e. old row 45 -> new row 3, old row 3 -> new row 4, ..., old row 44 -> new row 45.

$firstRow = 3
$lastRow = 45
$lastCol = "AY"

# Ensure text-like columns keep their text type (dates stored as text,
# collection numbers stored as text, etc.) instead of Excel auto-coercing
# them into numbers/dates when we write the rotated values back.
$textCols = @("C","D","F","G","H","I","P","T","U","V","W","Y","Z","AA","AB","AC","AR","AT","AW","AX","AY")
foreach ($col in $textCols) {
    $rng = $ws.Range($col + $firstRow + ":" + $col + $lastRow)
    $rng.NumberFormat = "@"
}

$fullRange = $ws.Range("A" + $firstRow + ":" + $lastCol + $lastRow)
$data = $fullRange.Value2

$nRows = $data.GetUpperBound(0)
$nCols = $data.GetUpperBound(1)

$newData = New-Object 'object[,]' $nRows,$nCols
for ($r = 1; $r -le $nRows; $r++) {
    $srcRow = $r - 1
    if ($srcRow -lt 1) { $srcRow = $nRows }
    for ($c = 1; $c -le $nCols; $c++) {
        $newData[$r-1, $c-1] = $data[$srcRow, $c]
    }
}

$fullRange.Value2 = $newData
